# Auto-generated edit script applying the Gilgamesh_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) on various rows
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 3664.743
$ws.Range("I15").Value = 3664.743
$ws.Range("K15").Value = 10994.229
$ws.Range("M15").Value = -10825.229
# Row 19
$ws.Range("H19").Value = 1065.6666
$ws.Range("I19").Value = 849.25
$ws.Range("K19").Value = 849.25
$ws.Range("M19").Value = -674.25
# Row 33
$ws.Range("H33").Value = 188
$ws.Range("I33").Value = 190.71428
$ws.Range("K33").Value = 190.71428
$ws.Range("M33").Value = 38.28572
# Row 113
$ws.Range("H113").Value = 3662.5
$ws.Range("I113").Value = 2650
$ws.Range("K113").Value = 2650
$ws.Range("M113").Value = 604
# Row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 834.75
$ws.Range("I2").Value = 580.0714
$ws.Range("J2").Value = 1191.3
$ws.Range("K2").Value = 580.0714
$ws.Range("L2").Value = 1191.3
$ws.Range("M2").Value = -467.0714
$ws.Range("N2").Value = -1417.3
# Row 7
$ws.Range("H7").Value = 49999
$ws.Range("J7").Value = 49999
$ws.Range("L7").Value = 49999
$ws.Range("N7").Value = -50227
# Row 32
$ws.Range("H32").Value = 4128.022
$ws.Range("I32").Value = 3047.8057
$ws.Range("J32").Value = 8016.8
$ws.Range("K32").Value = 3047.8057
$ws.Range("L32").Value = 8016.8
$ws.Range("M32").Value = -2760.8057
$ws.Range("N32").Value = -8590.799999999999
# Row 61
$ws.Range("H61").Value = 3770.25
$ws.Range("I61").Value = 2474.353
$ws.Range("J61").Value = 5773
$ws.Range("K61").Value = 2474.353
$ws.Range("L61").Value = 5773
$ws.Range("M61").Value = -2262.353
$ws.Range("N61").Value = -6197
# Row 107
$ws.Range("H107").Value = 32500
$ws.Range("J107").Value = 32500
$ws.Range("L107").Value = 32500
$ws.Range("N107").Value = -40180
# Row 116
$ws.Range("H116").Value = 834.75
$ws.Range("I116").Value = 580.0714
$ws.Range("J116").Value = 1191.3
$ws.Range("K116").Value = 580.0714
$ws.Range("L116").Value = 1191.3
$ws.Range("M116").Value = 1713.9286
$ws.Range("N116").Value = -5779.3

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 834.75
$ws.Range("I3").Value = 580.0714
$ws.Range("J3").Value = 1191.3
$ws.Range("K3").Value = 580.0714
$ws.Range("L3").Value = 1191.3
$ws.Range("M3").Value = -466.0714
$ws.Range("N3").Value = -1419.3
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 96
$ws.Range("H96").Value = 7499.5
$ws.Range("I96").Value = 7499.5
$ws.Range("K96").Value = 7499.5
$ws.Range("M96").Value = -4753.5
# Row 107
$ws.Range("H107").Value = 1418.5172
$ws.Range("I107").Value = 1214.4736
$ws.Range("J107").Value = 1806.2
$ws.Range("K107").Value = 1214.4736
$ws.Range("L107").Value = 1806.2
$ws.Range("M107").Value = 705.5264
$ws.Range("N107").Value = -5646.2
# Row 136
$ws.Range("H136").Value = 3770.25
$ws.Range("I136").Value = 2474.353
$ws.Range("J136").Value = 5773
$ws.Range("K136").Value = 7423.059
$ws.Range("L136").Value = 17319
$ws.Range("M136").Value = -4873.059
$ws.Range("N136").Value = -22419

$ws = $wb.Worksheets.Item("CRP")
# Row 63
$ws.Range("H63").Value = 95000
$ws.Range("J63").Value = 95000
$ws.Range("L63").Value = 95000
$ws.Range("N63").Value = -96372
# Row 66
$ws.Range("H66").Value = 95000
$ws.Range("J66").Value = 95000
$ws.Range("L66").Value = 285000
$ws.Range("N66").Value = -291864
# Row 99
$ws.Range("H99").Value = 6949.75
$ws.Range("I99").Value = 5000
$ws.Range("K99").Value = 5000
$ws.Range("M99").Value = -3502
# Row 126
$ws.Range("H126").Value = 6949.75
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530
# Row 132
$ws.Range("H132").Value = 3052.4814
$ws.Range("J132").Value = 3477
$ws.Range("L132").Value = 10431
$ws.Range("N132").Value = -15491
# Row 134
$ws.Range("H134").Value = 4193.1206
$ws.Range("I134").Value = 3971.8262
$ws.Range("J134").Value = 5041.4165
$ws.Range("K134").Value = 11915.4786
$ws.Range("L134").Value = 15124.2495
$ws.Range("M134").Value = -9380.4786
$ws.Range("N134").Value = -20194.2495

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 6995.9546
$ws.Range("J39").Value = 7239.1904
$ws.Range("L39").Value = 21717.5712
$ws.Range("N39").Value = -22305.5712
# Row 56
$ws.Range("H56").Value = 5289.3335
$ws.Range("I56").Value = 5289.3335
$ws.Range("K56").Value = 5289.3335
$ws.Range("M56").Value = -4759.3335
# Row 64
$ws.Range("H64").Value = 4562.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 4562.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 13687.5
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -14227.5
# Row 67
$ws.Range("H67").Value = 4562.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 4562.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 13687.5
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -15559.5
# Row 94
$ws.Range("H94").Value = 4921
$ws.Range("J94").Value = 5979.6
$ws.Range("L94").Value = 17938.8
$ws.Range("N94").Value = -19290.8
# Row 113
$ws.Range("H113").Value = 2376.923
$ws.Range("J113").Value = 2533.3333
$ws.Range("L113").Value = 7599.999899999999
$ws.Range("N113").Value = -11939.9999
# Row 125
$ws.Range("H125").Value = 4499.75
$ws.Range("J125").Value = 5333
$ws.Range("L125").Value = 15999
$ws.Range("N125").Value = -25839

$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 34999.5
$ws.Range("J52").Value = 34999.5
$ws.Range("L52").Value = 34999.5
$ws.Range("N52").Value = -35517.5
# Row 107
$ws.Range("H107").Value = 528.375
$ws.Range("I107").Value = 304.66666
$ws.Range("K107").Value = 304.66666
$ws.Range("M107").Value = 1615.33334
# Row 113
$ws.Range("H113").Value = 5208.0835
$ws.Range("I113").Value = 4900
$ws.Range("J113").Value = 5428.143
$ws.Range("K113").Value = 4900
$ws.Range("L113").Value = 5428.143
$ws.Range("M113").Value = -2730
$ws.Range("N113").Value = -9768.143
# Row 132
$ws.Range("H132").Value = 2427.1428
$ws.Range("I132").Value = 1220
$ws.Range("K132").Value = 3660
$ws.Range("M132").Value = -1130

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 338.72726
$ws.Range("I55").Value = 170.85715
$ws.Range("J55").Value = 632.5
$ws.Range("K55").Value = 170.85715
$ws.Range("L55").Value = 632.5
$ws.Range("M55").Value = 2.14285000000001
$ws.Range("N55").Value = -978.5
# Row 122
$ws.Range("H122").Value = 3430.4443
$ws.Range("I122").Value = 3666.3333
$ws.Range("J122").Value = 3312.5
$ws.Range("K122").Value = 10998.9999
$ws.Range("L122").Value = 9937.5
$ws.Range("M122").Value = -8548.999899999999
$ws.Range("N122").Value = -14837.5
# Row 132
$ws.Range("H132").Value = 16614.5
$ws.Range("I132").Value = 3488.8
$ws.Range("J132").Value = 25990
$ws.Range("K132").Value = 10466.4
$ws.Range("L132").Value = 77970
$ws.Range("M132").Value = -7936.400000000001
$ws.Range("N132").Value = -83030

$ws = $wb.Worksheets.Item("WVR")
# Row 52
$ws.Range("H52").Value = 19998.2
$ws.Range("I52").Value = 17499
$ws.Range("J52").Value = 29995
$ws.Range("K52").Value = 17499
$ws.Range("L52").Value = 29995
$ws.Range("M52").Value = -17273
$ws.Range("N52").Value = -30447
# Row 132
$ws.Range("H132").Value = 4232.517
$ws.Range("I132").Value = 4228.6924
$ws.Range("J132").Value = 4265.6665
$ws.Range("K132").Value = 12686.0772
$ws.Range("L132").Value = 12796.9995
$ws.Range("M132").Value = -10156.0772
$ws.Range("N132").Value = -17856.9995
# Row 136
$ws.Range("H136").Value = 5173.231
$ws.Range("I136").Value = 5717.3335
$ws.Range("J136").Value = 3949
$ws.Range("K136").Value = 17152.0005
$ws.Range("L136").Value = 11847
$ws.Range("M136").Value = -14602.0005
$ws.Range("N136").Value = -16947

